# Update timestamps for the "dc343afc-ddc9-4360-839a-25e47b7e753b" handback
# report row across the Overview, zh-cn, and de-de sheets, reflecting a
# regenerated report (later timestamps).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for dc343afc row (row 3)
$wsOverview.Range("G3").Value = "2016-09-02 02:52:45"

# zh-cn sheet: Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-09-02 02:52:41"
$wsZhCn.Range("K3").Value = "2016-09-02 02:53:00"

# de-de sheet: Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsDeDe.Range("H3").Value = "2016-09-02 02:52:45"
$wsDeDe.Range("K3").Value = "2016-09-02 02:53:14"
